$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1. Drop the old {contenidoHechos} paragraph plus the blank paragraph
#    that immediately follows it - the "Hechos" narrative now comes
#    exclusively from the {#hechos} repeating section below.
# ----------------------------------------------------------------------
$pContenido = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "{contenidoHechos}`r") {
        $pContenido = $p
        break
    }
}
if ($pContenido -ne $null) {
    $pAfter = $pContenido.Next()
    $d.Range($pContenido.Range.Start, $pAfter.Range.End).Delete() | Out-Null
}

# ----------------------------------------------------------------------
# 2. Turn the old single-image "imagenesHechos" loop into the new
#    per-fact "hechos" loop: {#hechos} ... {/hechos}
# ----------------------------------------------------------------------
$d.Content.Find.Execute("{#imagenesHechos}", $false, $false, $false, $false, $false, $true, 1, $false, "{#hechos}", 2) | Out-Null
$d.Content.Find.Execute("{/imagenesHechos}", $false, $false, $false, $false, $false, $true, 1, $false, "{/hechos}", 2) | Out-Null

# ----------------------------------------------------------------------
# 3. Rename {%src} -> {descripcionHecho} (the text of each fact)
# ----------------------------------------------------------------------
$d.Content.Find.Execute("{%src}", $false, $false, $false, $false, $false, $true, 1, $false, "{descripcionHecho}", 2) | Out-Null

# ----------------------------------------------------------------------
# 4. Add a new {%fotoHecho} paragraph (the image of each fact) right
#    after {descripcionHecho}, before the {/hechos} closing tag.
# ----------------------------------------------------------------------
$pDescripcion = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "{descripcionHecho}`r") {
        $pDescripcion = $p
        break
    }
}
if ($pDescripcion -ne $null) {
    $pDescripcion.Range.InsertParagraphAfter()
    $pFoto = $pDescripcion.Next()
    $pFoto.Range.Text = "{%fotoHecho}"
}

# ----------------------------------------------------------------------
# 5. Mark the footnote-related styles as Quick Styles (w:qFormat), as
#    they are now referenced by the new fact/image styling.
# ----------------------------------------------------------------------
$qFormatStyleNames = @("footnote text", "Unresolved Mention", "Texto nota pie Car")
foreach ($s in $d.Styles) {
    if ($qFormatStyleNames -contains $s.NameLocal) {
        $s.QuickStyle = $true
    }
}
